$d = $word.ActiveDocument
$sec = $d.Sections(1)

# footer2.xml (docPr id="2") is Word's Footers(1) -- PearsonLogo picture, image1.png -> image2.png
$sec.Footers(1).Range.InlineShapes(1).Name = "image2.png"

# footer1.xml (docPr id="3") is Word's Footers(2) -- PearsonLogo picture, image1.png -> image2.png
$sec.Footers(2).Range.InlineShapes(1).Name = "image2.png"

# header1.xml (docPr id="1") is Word's Headers(2) -- BTec_Logo-Orange picture, image2.jpg -> image1.jpg
$sec.Headers(2).Range.InlineShapes(1).Name = "image1.jpg"
